$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Add a header row to Sheet1: "Year" / "Jail Population"
$ws1.Range("A1").Value = "Year"
$ws1.Range("B1").Value = "Jail Population"

# Restore the selection on Sheet3 before moving away from it, so the
# cursor position sticks even though Sheet3 is no longer the active tab.
$ws3.Range("H14").Select()

# Finally, select Sheet1 (making it the active/visible tab) and set the
# zoom level + selected cell to match the reviewed state.
$ws1.Range("F8").Select()
$excel.ActiveWindow.Zoom = 185
